# New submission synced: 2026-02-08 22:15:27
# Target sheet: "JSS 3F" (form-response style sheet with Timestamp/Full Name/Admission No/AI Score)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3F")

# Row 4, column C ("Admission No" for Usman Muhammad Gubio) was entered as text "05";
# correct it to the true numeric value 5.
$ws.Range("C4").Value = 5

# Append the new form submission as row 5.
$ws.Range("A5").Value = "2026-02-08 22:15:27"
$ws.Range("B5").Value = "Usman Muhammad Gubio"

# Admission No "05" must stay text (leading zero preserved), same as the original C4 entry.
# Force text entry with a leading apostrophe, then reset the cell style so no extra
# number-format/style index gets attached to the cell.
$ws.Range("C5").Value = "'05"
$ws.Range("C5").Style = "Normal"

$ws.Range("D5").Value = 9
